$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching existing header style (copy format from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in I2:J16 values
$iValues = @(6, 5, 5, 7, 7, 6, 6, 5, 9, 7, 5, 5, 7, 2, 7)
$jValues = @(7, 6, 5, 8, 8, 6, 6, 6, 9, 8, 6, 7, 8, 3, 8)

for ($r = 0; $r -lt 15; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
